# Apply updates to column F (dSF) for specific rows as per repull of data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F12").Value = -3
$ws.Range("F16").Value = -4
$ws.Range("F18").Value = -1
$ws.Range("F24").Value = 3
$ws.Range("F26").Value = -2
$ws.Range("F35").Value = 3
$ws.Range("F37").Value = 0
